$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the header row (B1:G1) to the new, title-cased labels without the
# leading space, matching the "Add files via upload" / "Stats" section edit.
$ws.Range("B1").Value = "Pld"
$ws.Range("C1").Value = "Won"
$ws.Range("D1").Value = "Lost"
$ws.Range("E1").Value = "Tied"
$ws.Range("F1").Value = "Net RR"
$ws.Range("G1").Value = "Pts"
